# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 752-753) into the Uva price sheet,
# pushing the existing rows 752-778 down to 754-780.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 752.
$ws.Rows("752:753").Insert()

# --- New row 752: Flame Seedless, Provincia del Elquí ---
$ws.Range("A752").Value = 6
$ws.Range("B752").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C752").Value = "Metropolitana"
$ws.Range("D752").Value = 44568
$ws.Range("E752").Value = 13
$ws.Range("F752").Value = "Fruta"
$ws.Range("G752").Value = 100109
$ws.Range("H752").Value = "Uva"
$ws.Range("I752").Value = 100109001
$ws.Range("J752").Value = "Uva"
$ws.Range("K752").Value = "Flame Seedless"
$ws.Range("L752").Value = "Primera"
$ws.Range("M752").Value = 200
$ws.Range("N752").Value = 12000
$ws.Range("O752").Value = 12000
$ws.Range("P752").Value = 12000
$ws.Range("Q752").Value = "$/caja 15 kilos"
$ws.Range("R752").Value = "Provincia del Elquí"
$ws.Range("S752").Value = 800
$ws.Range("T752").Value = 15

# --- New row 753: Superior Seedless, Provincia del Elquí ---
$ws.Range("A753").Value = 6
$ws.Range("B753").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C753").Value = "Metropolitana"
$ws.Range("D753").Value = 44568
$ws.Range("E753").Value = 13
$ws.Range("F753").Value = "Fruta"
$ws.Range("G753").Value = 100109
$ws.Range("H753").Value = "Uva"
$ws.Range("I753").Value = 100109001
$ws.Range("J753").Value = "Uva"
$ws.Range("K753").Value = "Superior Seedless"
$ws.Range("L753").Value = "Primera"
$ws.Range("M753").Value = 200
$ws.Range("N753").Value = 15000
$ws.Range("O753").Value = 15000
$ws.Range("P753").Value = 15000
$ws.Range("Q753").Value = "$/caja 15 kilos"
$ws.Range("R753").Value = "Provincia del Elquí"
$ws.Range("S753").Value = 1000
$ws.Range("T753").Value = 15
